# Remove the trailing footer block consisting of the blank paragraph,
# the "Ver no Jupiter..." paragraph and the "© 2020 ..." paragraph that
# immediately follow the "Artigos de revistas..." bibliography line,
# leaving the remaining blank paragraph + page-break paragraph intact.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute(
    "Artigos de revistas especializadas e de jornais; Estudos, artigos, notícias e pesquisas via internet.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Resolve the paragraph index of the matched text so we can walk the
# paragraph collection reliably (relative Range navigation on the Find
# result was unreliable).
$anchorIndex = $d.Range(0, $rng.Start).Paragraphs.Count + 1

$blankPara = $d.Paragraphs.Item($anchorIndex + 1)
$copyrightPara = $d.Paragraphs.Item($anchorIndex + 3)

$deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
$deleteRange.Delete()
